{"js": "// Replace the 25 \"two-digit number divided by one-digit number\" problems\n// in the document's table cells, in document order. A couple of the\n// original problem strings repeat (e.g. \"46\u00f79=\" appears twice and maps to\n// two different replacements depending on position), so we search once per\n// distinct old string and then walk the matches in document order,\n// consuming the queued replacement values for that string one at a time.\n\n// Ordered list of [oldText, newText] exactly as they appear (top-to-bottom,\n// left-to-right) in the document body.\nconst replacements = [\n  [\"73\u00f78=\", \"53\u00f73=\"],\n  [\"97\u00f76=\", \"96\u00f74=\"],\n  [\"31\u00f78=\", \"30\u00f78=\"],\n  [\"88\u00f79=\", \"27\u00f72=\"],\n  [\"68\u00f75=\", \"59\u00f76=\"],\n  [\"78\u00f79=\", \"57\u00f78=\"],\n  [\"10\u00f79=\", \"59\u00f76=\"],\n  [\"46\u00f79=\", \"49\u00f78=\"],\n  [\"51\u00f79=\", \"37\u00f73=\"],\n  [\"66\u00f76=\", \"81\u00f75=\"],\n  [\"94\u00f72=\", \"61\u00f74=\"],\n  [\"44\u00f76=\", \"50\u00f78=\"],\n  [\"27\u00f77=\", \"42\u00f72=\"],\n  [\"50\u00f79=\", \"73\u00f79=\"],\n  [\"70\u00f73=\", \"58\u00f74=\"],\n  [\"27\u00f78=\", \"96\u00f77=\"],\n  [\"48\u00f79=\", \"44\u00f73=\"],\n  [\"39\u00f75=\", \"22\u00f73=\"],\n  [\"38\u00f74=\", \"47\u00f75=\"],\n  [\"53\u00f78=\", \"29\u00f79=\"],\n  [\"84\u00f75=\", \"89\u00f79=\"],\n  [\"59\u00f77=\", \"98\u00f78=\"],\n  [\"46\u00f79=\", \"57\u00f75=\"],\n  [\"90\u00f77=\", \"15\u00f77=\"],\n  [\"27\u00f74=\", \"63\u00f79=\"]\n];\n\n// Queue of intended new values per distinct old text, in the order they\n// must be applied (i.e. in document order of occurrence).\nconst queue = new Map();\nfor (const [oldText, newText] of replacements) {\n  if (!queue.has(oldText)) queue.set(oldText, []);\n  queue.get(oldText).push(newText);\n}\n\nconst body = context.document.body;\nconst searchResultsByOldText = new Map();\nfor (const oldText of queue.keys()) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items/text\");\n  searchResultsByOldText.set(oldText, results);\n}\n\nawait context.sync();\n\nfor (const [oldText, newTexts] of queue) {\n  const results = searchResultsByOldText.get(oldText);\n  if (results.items.length !== newTexts.length) {\n    throw new Error(\n      `Expected ${newTexts.length} occurrence(s) of \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newTexts[i], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"two-digit number divided by one-digit number\" problems\n# that live in the document's table cells, in document order (row-major,\n# left-to-right, top-to-bottom). A couple of the original problem strings\n# repeat (e.g. \"46\u00f79=\" appears twice and maps to two different\n# replacements depending on position), so we walk the table cells in\n# order and consume the ordered replacement list one entry at a time\n# rather than doing a single global find/replace.\n\n$d = $word.ActiveDocument\n\n# Ordered list of new values, matching the order the old values appear\n# top-to-bottom / left-to-right through the table.\n$newValues = @(\n  \"53\u00f73=\",\n  \"96\u00f74=\",\n  \"30\u00f78=\",\n  \"27\u00f72=\",\n  \"59\u00f76=\",\n  \"57\u00f78=\",\n  \"59\u00f76=\",\n  \"49\u00f78=\",\n  \"37\u00f73=\",\n  \"81\u00f75=\",\n  \"61\u00f74=\",\n  \"50\u00f78=\",\n  \"42\u00f72=\",\n  \"73\u00f79=\",\n  \"58\u00f74=\",\n  \"96\u00f77=\",\n  \"44\u00f73=\",\n  \"22\u00f73=\",\n  \"47\u00f75=\",\n  \"29\u00f79=\",\n  \"89\u00f79=\",\n  \"98\u00f78=\",\n  \"57\u00f75=\",\n  \"15\u00f77=\",\n  \"63\u00f79=\"\n)\n\n# Expected old values, purely as a sanity check so we fail loudly instead\n# of silently overwriting the wrong cell if the document doesn't look the\n# way we expect.\n$oldValues = @(\n  \"73\u00f78=\",\n  \"97\u00f76=\",\n  \"31\u00f78=\",\n  \"88\u00f79=\",\n  \"68\u00f75=\",\n  \"78\u00f79=\",\n  \"10\u00f79=\",\n  \"46\u00f79=\",\n  \"51\u00f79=\",\n  \"66\u00f76=\",\n  \"94\u00f72=\",\n  \"44\u00f76=\",\n  \"27\u00f77=\",\n  \"50\u00f79=\",\n  \"70\u00f73=\",\n  \"27\u00f78=\",\n  \"48\u00f79=\",\n  \"39\u00f75=\",\n  \"38\u00f74=\",\n  \"53\u00f78=\",\n  \"84\u00f75=\",\n  \"59\u00f77=\",\n  \"46\u00f79=\",\n  \"90\u00f77=\",\n  \"27\u00f74=\"\n)\n\n$t = $d.Tables.Item(1)\n$i = 0\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # A cell's Range.Text carries trailing control chars (paragraph mark\n    # \\r and/or the end-of-cell mark \\a) that aren't part of the visible\n    # content; strip them before comparing/counting.\n    $current = $cellRange.Text.TrimEnd([char]7, [char]13)\n\n    if ([string]::IsNullOrEmpty($current)) {\n      continue\n    }\n\n    if ($i -ge $oldValues.Count) {\n      continue\n    }\n\n    $expected = $oldValues[$i]\n    if ($current -ne $expected) {\n      throw \"Cell ($r,$c) text '$current' did not match expected '$expected' at position $i\"\n    }\n\n    $cellRange.Text = $newValues[$i]\n    $i = $i + 1\n  }\n}\n\nif ($i -ne $oldValues.Count) {\n  throw \"Only replaced $i of $($oldValues.Count) expected problems\"\n}\n"}
